$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a76c2a46d7b125e170bb96b094dbfb48074c1e35/e2e/1fbe89c7-79b7-42e9-be29-691beca46f43.md", "", "", "1fbe89c7-79b7-42e9-be29-691beca46f43.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/026284125cbe8635e937af3ed58cd29aa86f19dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/1fbe89c7-79b7-42e9-be29-691beca46f43.12050ce38d0e6915f2af4851e887c6d7e7120f54.zh-cn.xlf", "", "", "1fbe89c7-79b7-42e9-be29-691beca46f43.12050ce38d0e6915f2af4851e887c6d7e7120f54.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a76c2a46d7b125e170bb96b094dbfb48074c1e35/e2e/4ce294d5-985e-4915-bbeb-9ff69505b27e.md", "", "", "4ce294d5-985e-4915-bbeb-9ff69505b27e.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/026284125cbe8635e937af3ed58cd29aa86f19dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/4ce294d5-985e-4915-bbeb-9ff69505b27e.669fb09687659433288af44a0816314f06321646.zh-cn.xlf", "", "", "4ce294d5-985e-4915-bbeb-9ff69505b27e.669fb09687659433288af44a0816314f06321646.zh-cn.xlf")

$zh.Range("G2").Value = "2016-03-09 10:19:51"
$zh.Range("G3").Value = "2016-03-09 10:19:51"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a76c2a46d7b125e170bb96b094dbfb48074c1e35/e2e/1fbe89c7-79b7-42e9-be29-691beca46f43.md", "", "", "1fbe89c7-79b7-42e9-be29-691beca46f43.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d8556c08ae45616aced409f1fd40b67994f702f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/1fbe89c7-79b7-42e9-be29-691beca46f43.12050ce38d0e6915f2af4851e887c6d7e7120f54.de-de.xlf", "", "", "1fbe89c7-79b7-42e9-be29-691beca46f43.12050ce38d0e6915f2af4851e887c6d7e7120f54.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a76c2a46d7b125e170bb96b094dbfb48074c1e35/e2e/4ce294d5-985e-4915-bbeb-9ff69505b27e.md", "", "", "4ce294d5-985e-4915-bbeb-9ff69505b27e.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d8556c08ae45616aced409f1fd40b67994f702f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/4ce294d5-985e-4915-bbeb-9ff69505b27e.669fb09687659433288af44a0816314f06321646.de-de.xlf", "", "", "4ce294d5-985e-4915-bbeb-9ff69505b27e.669fb09687659433288af44a0816314f06321646.de-de.xlf")

$de.Range("G2").Value = "2016-03-09 10:19:58"
$de.Range("G3").Value = "2016-03-09 10:19:58"

Write-Host "Done applying handback updates"
